# SE-2339: Add simple instrument data.
# Adds a new "simple_instruments" worksheet (with shopping-centre sample
# rows) after the existing "term_deposits" sheet, then restores the
# "equities" sheet as the active/selected tab (matching the authored file).

$wb = $excel.ActiveWorkbook

# --- add the new worksheet at the end of the workbook -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "simple_instruments"

# --- header row ---------------------------------------------------------
$newSheet.Range("A1").Value = "Name"
$newSheet.Range("B1").Value = "ClientInternal"
$newSheet.Range("C1").Value = "Currency"
$newSheet.Range("D1").Value = "Class"
$newSheet.Range("E1").Value = "Type"

# --- first two instruments (names first) --------------------------------
$newSheet.Range("A2").Value = "LondonShoppingCentre"
$newSheet.Range("A3").Value = "BirminghamShoppingCentre"

$newSheet.Range("C2").Value = "GBP"
$newSheet.Range("D2").Value = "Unknown"
$newSheet.Range("E2").Value = "ShoppingCentres"

$newSheet.Range("C3").Value = "GBP"
$newSheet.Range("D3").Value = "Unknown"
$newSheet.Range("E3").Value = "ShoppingCentres"

# --- client-internal codes for all five instruments ----------------------
$newSheet.Range("B2").Value = "SHOPCENCI1"
$newSheet.Range("B3").Value = "SHOPCENCI2"
$newSheet.Range("B4").Value = "SHOPCENCI3"
$newSheet.Range("B5").Value = "SHOPCENCI4"
$newSheet.Range("B6").Value = "SHOPCENCI5"

# --- remaining three instruments ------------------------------------------
$newSheet.Range("A4").Value = "OxfordShoppingCentre"
$newSheet.Range("C4").Value = "GBP"
$newSheet.Range("D4").Value = "Unknown"
$newSheet.Range("E4").Value = "ShoppingCentres"

$newSheet.Range("A5").Value = "BathShoppingCentre"
$newSheet.Range("C5").Value = "GBP"
$newSheet.Range("D5").Value = "Unknown"
$newSheet.Range("E5").Value = "ShoppingCentres"

$newSheet.Range("A6").Value = "WarwickShoppingCentre"
$newSheet.Range("C6").Value = "GBP"
$newSheet.Range("D6").Value = "Unknown"
$newSheet.Range("E6").Value = "ShoppingCentres"

# leave the cursor on the new sheet parked away from the data, like the
# authored workbook (selection sits at I6, outside the A1:E6 table)
$newSheet.Range("I6").Select() | Out-Null

# --- restore "equities" as the active tab/selection -----------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate() | Out-Null
$ws1.Range("F4").Select() | Out-Null
